# "fundamental change of algorithms"
#  - Correct two stale "other" figures (I64, I130) from 564 to 0
#  - Update the sheet view: switch to right-to-left reading order, scroll
#    so row 124 is at the top, and leave I131 as the selected cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")
$ws.Activate()

# --- Data correction ---
$ws.Range("I64").Value = 0
$ws.Range("I130").Value = 0

# --- View settings ---
$window = $excel.ActiveWindow
$window.DisplayRightToLeft = $true
$window.ScrollRow = 124
$window.ScrollColumn = 1

$ws.Range("I131").Select()
